$wb = $excel.ActiveWorkbook

# The old ClassID value (a wasm/IBC-wrapped Juno address) is being replaced
# with a plain Juno address on the two sheets that referenced it ("A3" and
# "A5"). Updating the cell text automatically makes the shared-strings
# table drop the now-unused old string and append the new one, which is
# exactly the reindex shown in the diff.
$newValue = "juno1tzpn2jrz5mg7qrq32ceym7c8j7slc4m52zvyxnr0g3zefx3cmphsllzkc2"

$wsA3 = $wb.Worksheets.Item("A3")
$wsA3.Range("B2").Value = $newValue
$wsA3.Activate()
$wsA3.Range("B3").Select()

$wsA5 = $wb.Worksheets.Item("A5")
$wsA5.Range("B2").Value = $newValue
$wsA5.Activate()
$wsA5.Range("C5").Select()
